$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = [System.Convert]::ToDouble("0.07316220184780846")
$ws.Range("C2").Value = [System.Convert]::ToDouble("6.893845046377231e-11")
$ws.Range("D2").Value = [System.Convert]::ToDouble("3.785670380641082e-12")
$ws.Range("B3").Value = [System.Convert]::ToDouble("0.1492751576917354")
$ws.Range("C3").Value = [System.Convert]::ToDouble("1.641842548698402e-11")
$ws.Range("D3").Value = [System.Convert]::ToDouble("1.608402642568462e-12")
$ws.Range("B4").Value = [System.Convert]::ToDouble("0.1342327798802501")
$ws.Range("C4").Value = [System.Convert]::ToDouble("7.098803971480968e-12")
$ws.Range("D4").Value = [System.Convert]::ToDouble("7.532708734793558e-13")
$ws.Range("B5").Value = [System.Convert]::ToDouble("0.09516638299361919")
$ws.Range("C5").Value = [System.Convert]::ToDouble("3.999144734425883e-12")
$ws.Range("D5").Value = [System.Convert]::ToDouble("5.593635506137449e-13")
$ws.Range("B6").Value = [System.Convert]::ToDouble("0.1293556144447252")
$ws.Range("C6").Value = [System.Convert]::ToDouble("2.608919824226776e-12")
$ws.Range("D6").Value = [System.Convert]::ToDouble("3.576458112233485e-13")
$ws.Range("B7").Value = [System.Convert]::ToDouble("0.06840164815886485")
$ws.Range("C7").Value = [System.Convert]::ToDouble("1.800924621303047e-12")
$ws.Range("D7").Value = [System.Convert]::ToDouble("9.107970217570611e-14")
$ws.Range("B8").Value = [System.Convert]::ToDouble("0.1202209034306793")
$ws.Range("C8").Value = [System.Convert]::ToDouble("1.312258098878103e-12")
$ws.Range("D8").Value = [System.Convert]::ToDouble("1.471198072053214e-13")
$ws.Range("B9").Value = [System.Convert]::ToDouble("0.07495653084902071")
$ws.Range("C9").Value = [System.Convert]::ToDouble("1.061983234615559e-12")
$ws.Range("D9").Value = [System.Convert]::ToDouble("8.300249088131099e-14")
$ws.Range("B10").Value = [System.Convert]::ToDouble("0.1034943423535781")
$ws.Range("C10").Value = [System.Convert]::ToDouble("8.283359730879508e-13")
$ws.Range("D10").Value = [System.Convert]::ToDouble("1.193411934089963e-13")
$ws.Range("B11").Value = [System.Convert]::ToDouble("0.12282583880371")
$ws.Range("C11").Value = [System.Convert]::ToDouble("6.909482715916617e-13")
$ws.Range("D11").Value = [System.Convert]::ToDouble("6.525583682174036e-14")
